$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the time columns (B:D) keep the existing "h:mm" time-text number format
# for the new rows (style index 1 in the original workbook).
$ws.Range("B2:D5").NumberFormat = "h:mm"

# --- Row 1 (header/first data row) - timestamps updated, column A unchanged ---
$ws.Range("A1").Value = "AS"
$ws.Range("B1").Value = "2025-02-03T11:30"
$ws.Range("C1").Value = "2025-02-03T13:00"
$ws.Range("D1").Value = "2025-02-03T14:50"

# --- Row 2 ---
$ws.Range("A2").Value = "AS"
$ws.Range("B2").Value = "2025-02-04T14:24"
$ws.Range("C2").Value = "2025-02-04T16:25"
$ws.Range("D2").Value = "2025-02-04T17:25"
$ws.Range("E2").Value = "rgb(30, 144, 255)"

# --- Row 3 ---
$ws.Range("A3").Value = "AS"
$ws.Range("B3").Value = "2025-02-04T15:24"
$ws.Range("C3").Value = "2025-02-04T17:25"
$ws.Range("D3").Value = "2025-02-04T18:25"
$ws.Range("E3").Value = "rgb(46, 139, 87)"

# --- Row 4 (new) ---
$ws.Range("A4").Value = "AS"
$ws.Range("B4").Value = "2025-02-06T12:28"
$ws.Range("C4").Value = "2025-02-04T14:28"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "rgb(194, 24, 7)"

# --- Row 5 (new) ---
$ws.Range("A5").Value = "AS"
$ws.Range("B5").Value = "2025-02-06T11:28"
$ws.Range("C5").Value = "2025-02-04T15:28"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "rgb(241, 156, 187)"
